# Update design layout and add new table style
# Rewrites the data rows of the "yesterday no sales" report:
#  - refreshes the BSL/ISL numbers and UOM values for existing brands
#  - drops "Sk-Mox" and introduces two new brands: "Quinox" and "Volmax"
#  - appends a new row for "Zithrox 250mg Tablet - 6's"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = 4
$ws.Cells.Item(2,2).Value = "Desodin"
$ws.Cells.Item(2,3).Value = 9
$ws.Cells.Item(2,4).Value = "Desodin 60ml Syrup"
$ws.Cells.Item(2,5).Value = "60 ml"
$ws.Cells.Item(3,1).Value = 5
$ws.Cells.Item(3,2).Value = "Dinafex"
$ws.Cells.Item(3,3).Value = 10
$ws.Cells.Item(3,4).Value = "Dinafex 60mg Tablet"
$ws.Cells.Item(3,5).Value = "30's"
$ws.Cells.Item(4,1).Value = 5
$ws.Cells.Item(4,2).Value = "Dinafex"
$ws.Cells.Item(4,3).Value = 11
$ws.Cells.Item(4,4).Value = "Dinafex 180mg Tablet"
$ws.Cells.Item(4,5).Value = "30's"
$ws.Cells.Item(5,1).Value = 5
$ws.Cells.Item(5,2).Value = "Dinafex"
$ws.Cells.Item(5,3).Value = 15
$ws.Cells.Item(5,4).Value = "Dinafex 120mg Tablet"
$ws.Cells.Item(5,5).Value = "30's"
$ws.Cells.Item(6,1).Value = 6
$ws.Cells.Item(6,2).Value = "Dorenta"
$ws.Cells.Item(6,3).Value = 18
$ws.Cells.Item(6,4).Value = "Dorenta 50mg Tablet"
$ws.Cells.Item(6,5).Value = "50's"
$ws.Cells.Item(7,1).Value = 7
$ws.Cells.Item(7,2).Value = "Etorix"
$ws.Cells.Item(7,3).Value = 19
$ws.Cells.Item(7,4).Value = "Etorix 120mg Tablet"
$ws.Cells.Item(7,5).Value = "20's"
$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = "Etorix"
$ws.Cells.Item(8,3).Value = 23
$ws.Cells.Item(8,4).Value = "Etorix 90mg Tablet"
$ws.Cells.Item(8,5).Value = "30's"
$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "Etorix"
$ws.Cells.Item(9,3).Value = 24
$ws.Cells.Item(9,4).Value = "Etorix 60mg Tablet - 40's"
$ws.Cells.Item(9,5).Value = "40's"
$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "Fenobac"
$ws.Cells.Item(10,3).Value = 26
$ws.Cells.Item(10,4).Value = "Fenobac 100ml Syrup"
$ws.Cells.Item(10,5).Value = "100ml"
$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "Flucloxin"
$ws.Cells.Item(11,3).Value = 28
$ws.Cells.Item(11,4).Value = "Flucloxin 500mg Capsule - 36's"
$ws.Cells.Item(11,5).Value = "36 's"
$ws.Cells.Item(12,1).Value = 9
$ws.Cells.Item(12,2).Value = "Flucloxin"
$ws.Cells.Item(12,3).Value = 29
$ws.Cells.Item(12,4).Value = "Flucloxin 500mg Capsule"
$ws.Cells.Item(12,5).Value = "30 's"
$ws.Cells.Item(13,1).Value = 10
$ws.Cells.Item(13,2).Value = "Geminox"
$ws.Cells.Item(13,3).Value = 34
$ws.Cells.Item(13,4).Value = "Geminox 320mg Tablet - 8's"
$ws.Cells.Item(13,5).Value = "8 's"
$ws.Cells.Item(14,1).Value = 11
$ws.Cells.Item(14,2).Value = "Ketonic"
$ws.Cells.Item(14,3).Value = 35
$ws.Cells.Item(14,4).Value = "Ketonic 30mg IM/IV Injection - 4's"
$ws.Cells.Item(14,5).Value = "4's"
$ws.Cells.Item(15,1).Value = 11
$ws.Cells.Item(15,2).Value = "Ketonic"
$ws.Cells.Item(15,3).Value = 36
$ws.Cells.Item(15,4).Value = "Ketonic 30mg Injection"
$ws.Cells.Item(15,5).Value = "5 's"
$ws.Cells.Item(16,1).Value = 11
$ws.Cells.Item(16,2).Value = "Ketonic"
$ws.Cells.Item(16,3).Value = 40
$ws.Cells.Item(16,4).Value = "Ketonic 10mg Tablet"
$ws.Cells.Item(16,5).Value = "20's"
$ws.Cells.Item(17,1).Value = 12
$ws.Cells.Item(17,2).Value = "Kynol"
$ws.Cells.Item(17,3).Value = 41
$ws.Cells.Item(17,4).Value = "Kynol D 25mg Tablet"
$ws.Cells.Item(17,5).Value = "60 's"
$ws.Cells.Item(18,1).Value = 12
$ws.Cells.Item(18,2).Value = "Kynol"
$ws.Cells.Item(18,3).Value = 42
$ws.Cells.Item(18,4).Value = "Kynol TR 100mg Capsule"
$ws.Cells.Item(18,5).Value = "50 's"
$ws.Cells.Item(19,1).Value = 12
$ws.Cells.Item(19,2).Value = "Kynol"
$ws.Cells.Item(19,3).Value = 43
$ws.Cells.Item(19,4).Value = "Kynol TR 200mg Capsule"
$ws.Cells.Item(19,5).Value = "30 's"
$ws.Cells.Item(20,1).Value = 17
$ws.Cells.Item(20,2).Value = "Naprox"
$ws.Cells.Item(20,3).Value = 58
$ws.Cells.Item(20,4).Value = "Naprox Plus 500mg Tablet - 30's"
$ws.Cells.Item(20,5).Value = "30 's"
$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "Oradin"
$ws.Cells.Item(21,3).Value = 64
$ws.Cells.Item(21,4).Value = "Oradin Plus Tablet - 40's"
$ws.Cells.Item(21,5).Value = "40 's"
$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "Osticare"
$ws.Cells.Item(22,3).Value = 68
$ws.Cells.Item(22,4).Value = "Osticare Tablet 24's"
$ws.Cells.Item(22,5).Value = "24's"
$ws.Cells.Item(23,1).Value = 22
$ws.Cells.Item(23,2).Value = "Quinox"
$ws.Cells.Item(23,3).Value = 73
$ws.Cells.Item(23,4).Value = "Quinox 100ml IV Infusion"
$ws.Cells.Item(23,5).Value = "1's"
$ws.Cells.Item(24,1).Value = 23
$ws.Cells.Item(24,2).Value = "Rupaday"
$ws.Cells.Item(24,3).Value = 77
$ws.Cells.Item(24,4).Value = "Rupaday Oral Solution 60ml"
$ws.Cells.Item(24,5).Value = "1's"
$ws.Cells.Item(25,1).Value = 33
$ws.Cells.Item(25,2).Value = "Volmax"
$ws.Cells.Item(25,3).Value = 96
$ws.Cells.Item(25,4).Value = "Volmax SR 100mg Capsule"
$ws.Cells.Item(25,5).Value = "60 's"
$ws.Cells.Item(26,1).Value = 35
$ws.Cells.Item(26,2).Value = "Zithrox"
$ws.Cells.Item(26,3).Value = 100
$ws.Cells.Item(26,4).Value = "Zithrox 500mg Tablet"
$ws.Cells.Item(26,5).Value = "6 's"
$ws.Cells.Item(27,1).Value = 35
$ws.Cells.Item(27,2).Value = "Zithrox"
$ws.Cells.Item(27,3).Value = 102
$ws.Cells.Item(27,4).Value = "Zithrox 15ml Suspension"
$ws.Cells.Item(27,5).Value = "15 ml"
$ws.Cells.Item(28,1).Value = 35
$ws.Cells.Item(28,2).Value = "Zithrox"
$ws.Cells.Item(28,3).Value = 104
$ws.Cells.Item(28,4).Value = "Zithrox 30ml Dry Suspension"
$ws.Cells.Item(28,5).Value = "30ml"
$ws.Cells.Item(29,1).Value = 35
$ws.Cells.Item(29,2).Value = "Zithrox"
$ws.Cells.Item(29,3).Value = 105
$ws.Cells.Item(29,4).Value = "Zithrox 250mg Tablet - 6's"
$ws.Cells.Item(29,5).Value = "6's"
